$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row additions (P1, Q1) - styled like other header cells (bold + border)
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

$hdr = $ws.Range("P1:Q1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108  # xlCenter
$hdr.VerticalAlignment = -4160   # xlTop
$hdr.Borders.LineStyle = 1       # xlContinuous
$hdr.Borders.Weight = 2          # xlThin

# Update existing columns I, K, M, O for rows 2-25
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2   # I -> 2
    $ws.Cells.Item($r, 11).Value = 1  # K -> 1
    $ws.Cells.Item($r, 13).Value = 2  # M -> 2
    $ws.Cells.Item($r, 15).Value = 1  # O -> 1
    $ws.Cells.Item($r, 16).Value = 2  # P -> 2 (new)
    $ws.Cells.Item($r, 17).Value = 2  # Q -> 2 (new)
}
